$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.188.35"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").Value = "1.912.09"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "'0.8202"
$ws.Range("E5").Value = "  +4.36%  "

$ws.Range("D6").Value = "'243.74"
$ws.Range("E6").Value = "  +0.49%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").Value = "'0.3253"
$ws.Range("E8").Value = "  +3.31%  "

$ws.Range("D9").Value = "'26.83"
$ws.Range("E9").Value = "  +2.72%  "

$ws.Range("D10").Value = "'0.07089"
$ws.Range("E10").Value = "  +2.85%  "

$ws.Range("D11").Value = "'0.08081"
$ws.Range("E11").Value = "  +1.38%  "

$ws.Range("D12").Value = "'0.7784"
$ws.Range("E12").Value = "  +4.77%  "

$ws.Range("D13").Value = "1.902.75"
$ws.Range("E13").Value = "  -0.37%  "

$ws.Range("D14").Value = "'5.347"
$ws.Range("E14").Value = "  +2.66%  "

$ws.Range("D15").Value = "'93.41"
$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("D16").Value = "30.204.52"
$ws.Range("E16").Value = "  +0.55%  "

$ws.Range("E17").Value = "  +2.32%  "

$ws.Range("D18").Value = "'5.943"
$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("D19").Value = "'246.72"
$ws.Range("E19").Value = "  +0.30%  "

$ws.Range("D20").Value = "'0.000007819"
$ws.Range("E20").Value = "  +0.91%  "

$ws.Range("D21").Value = "2.163.10"
$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("D23").Value = "'0.9997"
$ws.Range("E23").Value = "  -0.29%  "

$ws.Range("D24").Value = "'7.430"
$ws.Range("E24").Value = "  +8.11%  "

$ws.Range("D25").Value = "'0.1681"
$ws.Range("E25").Value = "  +23.10%  "

$ws.Range("D26").Value = "'9.387"
$ws.Range("E26").Value = "  +1.24%  "

$ws.Range("D27").Value = "'167.78"
$ws.Range("E27").Value = "  -1.20%  "

$ws.Range("D28").Value = "'19.01"
$ws.Range("E28").Value = "  +0.60%  "

$ws.Range("D29").Value = "'2.114"
$ws.Range("E29").Value = "  +4.25%  "

$ws.Range("D30").Value = "'1.377"
$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("D31").Value = "'1.531"
$ws.Range("E31").Value = "  +0.67%  "

$ws.Range("D32").Value = "'4.313"
$ws.Range("E32").Value = "  -0.31%  "

$ws.Range("D33").Value = "'0.05748"
$ws.Range("E33").Value = "  +5.55%  "

$ws.Range("D34").Value = "'4.107"
$ws.Range("E34").Value = "  +0.14%  "

$ws.Range("D35").Value = "'1.276"
$ws.Range("E35").Value = "  +1.76%  "

$ws.Range("D36").Value = "'0.7391"
$ws.Range("E36").Value = "  +0.70%  "

$ws.Range("E37").Value = "  +0.25%  "

$ws.Range("D38").Value = "'2.719"
$ws.Range("E38").Value = "  -0.29%  "

$ws.Range("D39").Value = "'0.01930"
$ws.Range("E39").Value = "  -0.22%  "

$ws.Range("E40").Value = "  +0.35%  "

$ws.Range("D41").Value = "'0.4470"
$ws.Range("E41").Value = "  +1.06%  "

$ws.Range("D42").Value = "'73.57"
$ws.Range("E42").Value = "  +2.07%  "

$ws.Range("D43").Value = "'5.966"
$ws.Range("E43").Value = "  -2.94%  "

$ws.Range("D44").Value = "'0.8502"
$ws.Range("E44").Value = "  +1.75%  "

$ws.Range("D45").Value = "'1.918"
$ws.Range("E45").Value = "  +1.97%  "

$ws.Range("D46").Value = "'0.9996"
$ws.Range("E46").Value = "  -0.23%  "

$ws.Range("D47").Value = "'103.07"
$ws.Range("E47").Value = "  +2.62%  "

$ws.Range("D48").Value = "1.025.15"
$ws.Range("E48").Value = "  +5.01%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.883"
$ws.Range("E49").Value = "  +0.78%  "

$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "'7.597"
$ws.Range("E50").Value = "  +1.12%  "

$ws.Range("D51").Value = "'1.574"
$ws.Range("E51").Value = "  +5.65%  "
